# Add daily COVID bulletin rows for data up to January 25th (rows 184-208)
# to "Planilha1" (sheet1), then leave that sheet active with A208 selected,
# matching the scroll/freeze-pane state saved by the author.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")

# New rows: date-serial, then columns B..J as recorded in the bulletin.
$newRows = @(
    @(44197,4445,47,874,5366,807,48,2,46,19),
    @(44198,4466,32,877,5375,811,47,2,45,19),
    @(44199,4467,31,878,5376,814,45,2,43,19),
    @(44200,4467,43,878,5388,816,41,2,39,21),
    @(44201,4471,113,879,5463,818,40,2,38,21),
    @(44202,4471,134,882,5487,822,39,2,37,21),
    @(44203,4527,103,902,5532,828,52,2,50,22),
    @(44204,4537,138,905,5580,831,52,2,50,22),
    @(44205,4591,115,932,5638,838,72,3,69,22),
    @(44206,4591,120,932,5643,843,67,4,63,22),
    @(44207,4608,137,938,5683,860,56,5,51,22),
    @(44208,4616,172,941,5729,865,54,5,49,22),
    @(44209,4660,160,961,5781,870,69,5,64,22),
    @(44210,4733,120,973,5826,876,75,5,70,22),
    @(44211,4768,88,993,5849,884,86,5,81,23),
    @(44212,4789,92,1002,5883,886,93,5,88,23),
    @(44213,4792,98,1002,5892,889,90,5,85,23),
    @(44214,4824,81,1011,5916,900,88,3,85,23),
    @(44215,4837,120,1022,5979,912,87,3,84,23),
    @(44216,4868,120,1034,6022,925,86,4,82,23),
    @(44217,4903,112,1050,6065,932,95,4,91,23),
    @(44218,4937,104,1068,6109,946,99,4,95,23),
    @(44219,4991,93,1080,6164,952,105,5,100,23),
    @(44220,5014,62,1090,6166,957,110,5,105,23),
    @(44221,5017,93,1091,6201,968,100,5,95,23)
)

$startRow = 184
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws1.Cells.Item($r, $col).Value = $rowValues[$col - 1]
    }
}

# Leave Planilha1 as the active sheet/selection, matching the saved view.
$ws1.Activate()
$ws1.Range("A208").Select()
